$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Tgfb2", "Tgfbr2", "ECs", 2, 0.6666666666666666, 1.465615333333333, 4.396846, 0.04672291954663727, 0.04672291954663728, 3, 1, 30.84053233333333, 92.521597, 0.2792762306509579, 0.2792762306509579, 45.20035707589577, 406.803213683062, 0.01304860085599282, 0.01304860085599282),
    @("ECs", "Tgfb2", "Tgfbr2", "FAPs", 2, 0.6666666666666666, 1.465615333333333, 4.396846, 0.04672291954663727, 0.04672291954663728, 3, 1, 69.52746833333333, 208.582405, 0.6296055163046042, 0.6296055163046042, 101.9005236771811, 917.10471309463, 0.02941700788441904, 0.02941700788441905),
    @("ECs", "Tgfb2", "Tgfbr2", "sCs", 2, 0.6666666666666666, 1.465615333333333, 4.396846, 0.04672291954663727, 0.04672291954663728, 3, 1, 10.06220766666667, 30.186623, 0.09111825304443795, 0.09111825304443795, 14.74732584345089, 132.725932591058, 0.004257310806225411, 0.004257310806225411),
    @("FAPs", "Tgfb2", "Tgfbr2", "ECs", 3, 1, 18.88237266666667, 56.64711800000001, 0.6019584804341267, 0.6019584804341268, 3, 1, 30.84053233333333, 92.521597, 0.2792762306509579, 0.2792762306509579, 582.3424247563829, 5241.081822807447, 0.1681126954240213, 0.1681126954240213),
    @("FAPs", "Tgfb2", "Tgfbr2", "FAPs", 3, 1, 18.88237266666667, 56.64711800000001, 0.6019584804341267, 0.6019584804341268, 3, 1, 69.52746833333333, 208.582405, 0.6296055163046042, 0.6296055163046042, 1312.843567639866, 11815.59210875879, 0.3789963798676633, 0.3789963798676634),
    @("FAPs", "Tgfb2", "Tgfbr2", "sCs", 3, 1, 18.88237266666667, 56.64711800000001, 0.6019584804341267, 0.6019584804341268, 3, 1, 10.06220766666667, 30.186623, 0.09111825304443795, 0.09111825304443795, 189.9983550113905, 1709.985195102514, 0.0548494051424421, 0.05484940514244212),
    @("sCs", "Tgfb2", "Tgfbr2", "ECs", 3, 1, 11.020243, 33.060729, 0.351318600019236, 0.351318600019236, 3, 1, 30.84053233333333, 92.521597, 0.2792762306509579, 0.2792762306509579, 339.8701605626903, 3058.831445064213, 0.09811493437094379, 0.09811493437094379),
    @("sCs", "Tgfb2", "Tgfbr2", "FAPs", 3, 1, 11.020243, 33.060729, 0.351318600019236, 0.351318600019236, 3, 1, 69.52746833333333, 208.582405, 0.6296055163046042, 0.6296055163046042, 766.2095962081384, 6895.886365873245, 0.2211921285525218, 0.2211921285525218),
    @("sCs", "Tgfb2", "Tgfbr2", "sCs", 3, 1, 11.020243, 33.060729, 0.351318600019236, 0.351318600019236, 3, 1, 10.06220766666667, 30.186623, 0.09111825304443795, 0.09111825304443795, 110.8879736031297, 997.991762428167, 0.03201153709577043, 0.03201153709577043)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    for ($j = 0; $j -lt $data[$i].Count; $j++) {
        $c = $j + 1
        $ws.Cells.Item($r, $c).Value = $data[$i][$j]
    }
}